# Apply updated dSF (column F) values to Sheet1, rows 2-18.
# This reflects a data repull / recalculated mean for the dSF column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    3  = -6
    4  = -8
    6  = -5
    7  = 2
    8  = -6
    9  = -3
    10 = 2
    11 = 1
    12 = -4
    13 = -2
    14 = -9
    15 = 1
    16 = -5
    17 = 9
    18 = 2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
